$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename row-1 header cells from the old human-readable labels to the new
# snake_case variable-style names (columns M:X, the P1..P12 headers, are
# left untouched).
$ws.Range("B1").Value = "n_customers"
$ws.Range("C1").Value = "n_vehicles"
$ws.Range("D1").Value = "depot_location"
$ws.Range("E1").Value = "dispersion_distance_depot"
$ws.Range("F1").Value = "dispersion_customers"
$ws.Range("G1").Value = "perc_pattern_school"
$ws.Range("H1").Value = "perc_pattern_home"
$ws.Range("I1").Value = "perc_pattern_restaurant"
$ws.Range("J1").Value = "expected_demand_on_tot_capacity"
$ws.Range("K1").Value = "initial_demand_on_tot_capacity"
$ws.Range("L1").Value = "time_slot_size"

# Row 2 previously carried example/placeholder values under DEPOT LOCATION,
# the two "deviation distances..." columns and the TIMESLOT SIZE column.
# Drop D2:F2 entirely and just clear L2's content (keeping its centered
# style) so the sheet is ready for fresh training data entry.
$ws.Range("D2:F2").Clear()
$ws.Range("L2").ClearContents()
